$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data cell by cell.
# NumberFormat "@" (Text) is applied before assignment and cleared
# afterwards so that numeric-looking strings (e.g. "1.003") are stored
# as text (matching the original inline/shared string cell type)
# instead of being auto-converted to numbers by Excel, while leaving
# the cell style untouched (ClearFormats drops the temporary format).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "28.269.69"
Set-TextValue $ws.Range("E2") "  +0.97%  "
Set-TextValue $ws.Range("D3") "1.805.56"
Set-TextValue $ws.Range("E3") "  +2.84%  "
Set-TextValue $ws.Range("D4") "1.003"
Set-TextValue $ws.Range("E4") "  -0.03%  "
Set-TextValue $ws.Range("D5") "338.31"
Set-TextValue $ws.Range("E5") "  +0.80%  "
Set-TextValue $ws.Range("D6") "0.9994"
Set-TextValue $ws.Range("E6") "  +0.02%  "
Set-TextValue $ws.Range("D7") "0.4707"
Set-TextValue $ws.Range("E7") "  +22.72%  "
Set-TextValue $ws.Range("D8") "0.3805"
Set-TextValue $ws.Range("E8") "  +11.51%  "
Set-TextValue $ws.Range("D9") "45.35"
Set-TextValue $ws.Range("E9") "  -1.06%  "
Set-TextValue $ws.Range("D10") "1.152"
Set-TextValue $ws.Range("D11") "0.07620"
Set-TextValue $ws.Range("E11") "  +5.19%  "
Set-TextValue $ws.Range("B12") "Solana"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D12") "22.43"
Set-TextValue $ws.Range("E12") "  -0.81%  "
Set-TextValue $ws.Range("B13") "BinanceUSD"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D13") "1.001"
Set-TextValue $ws.Range("E13") "  -0.01%  "
Set-TextValue $ws.Range("D14") "6.338"
Set-TextValue $ws.Range("D15") "7.458"
Set-TextValue $ws.Range("E15") "  +4.28%  "
Set-TextValue $ws.Range("D16") "1.808.37"
Set-TextValue $ws.Range("E16") "  +3.37%  "
Set-TextValue $ws.Range("D17") "0.00001094"
Set-TextValue $ws.Range("E17") "  +2.83%  "
Set-TextValue $ws.Range("D18") "0.06720"
Set-TextValue $ws.Range("E18") "  +1.71%  "
Set-TextValue $ws.Range("D19") "81.87"
Set-TextValue $ws.Range("E19") "  +3.09%  "
Set-TextValue $ws.Range("D20") "0.9994"
Set-TextValue $ws.Range("E20") "  +0.02%  "
Set-TextValue $ws.Range("D21") "17.44"
Set-TextValue $ws.Range("E21") "  +3.84%  "
Set-TextValue $ws.Range("D22") "6.412"
Set-TextValue $ws.Range("E22") "  +3.31%  "
Set-TextValue $ws.Range("D23") "28.264.57"
Set-TextValue $ws.Range("E23") "  +0.93%  "
Set-TextValue $ws.Range("D24") "11.87"
Set-TextValue $ws.Range("E24") "  +1.52%  "
Set-TextValue $ws.Range("D25") "2.407"
Set-TextValue $ws.Range("E25") "  +1.14%  "
Set-TextValue $ws.Range("E26") "  +4.15%  "
Set-TextValue $ws.Range("D27") "153.81"
Set-TextValue $ws.Range("E27") "  -0.19%  "
Set-TextValue $ws.Range("D28") "2.371"
Set-TextValue $ws.Range("E28") "  +2.73%  "
Set-TextValue $ws.Range("D29") "2.008.61"
Set-TextValue $ws.Range("E29") "  +3.00%  "
Set-TextValue $ws.Range("D30") "133.16"
Set-TextValue $ws.Range("E30") "  +1.52%  "
Set-TextValue $ws.Range("D31") "1.256"
Set-TextValue $ws.Range("E31") "  -0.57%  "
Set-TextValue $ws.Range("D32") "4.035"
Set-TextValue $ws.Range("E32") "  +0.16%  "
Set-TextValue $ws.Range("D33") "0.09638"
Set-TextValue $ws.Range("E33") "  +9.19%  "
Set-TextValue $ws.Range("D34") "5.866"
Set-TextValue $ws.Range("E34") "  -0.01%  "
Set-TextValue $ws.Range("D35") "0.2248"
Set-TextValue $ws.Range("E35") "  +6.60%  "
Set-TextValue $ws.Range("B36") "Aptos"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D36") "12.13"
Set-TextValue $ws.Range("E36") "  -0.92%  "
Set-TextValue $ws.Range("B37") "Hedera"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D37") "0.06373"
Set-TextValue $ws.Range("E37") "  +3.38%  "
Set-TextValue $ws.Range("D38") "0.02357"
Set-TextValue $ws.Range("E38") "  +2.89%  "
Set-TextValue $ws.Range("D39") "5.254"
Set-TextValue $ws.Range("E39") "  +1.80%  "
Set-TextValue $ws.Range("D40") "0.6632"
Set-TextValue $ws.Range("E40") "  +0.61%  "
Set-TextValue $ws.Range("B41") "WEMIXTOKEN"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D41") "1.503"
Set-TextValue $ws.Range("E41") "  -2.94%  "
Set-TextValue $ws.Range("B42") "TrustWalletToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D42") "1.237"
Set-TextValue $ws.Range("E42") "  +1.63%  "
Set-TextValue $ws.Range("D43") "8.247"
Set-TextValue $ws.Range("E43") "  +3.41%  "
Set-TextValue $ws.Range("D44") "14.14"
Set-TextValue $ws.Range("E44") "  +3.05%  "
Set-TextValue $ws.Range("D45") "0.9991"
Set-TextValue $ws.Range("E45") "  +0.02%  "
Set-TextValue $ws.Range("D46") "0.6136"
Set-TextValue $ws.Range("E46") "  +1.23%  "
Set-TextValue $ws.Range("D47") "3.852"
Set-TextValue $ws.Range("E47") "  +0.39%  "
Set-TextValue $ws.Range("D48") "130.37"
Set-TextValue $ws.Range("E48") "  +2.42%  "
Set-TextValue $ws.Range("D49") "2.034"
Set-TextValue $ws.Range("E49") "  +1.19%  "
Set-TextValue $ws.Range("D50") "0.07157"
Set-TextValue $ws.Range("E50") "  +2.55%  "
Set-TextValue $ws.Range("D51") "1.178"
Set-TextValue $ws.Range("E51") "  +0.69%  "
